$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("HV0X" -> "HV00X"), set column-by-column (all of D, then all of E)
# so the rebuilt shared-strings table groups HV-entries before HP-entries,
# matching how Excel coalesces the string table after these edits.
$ws.Range("D1").Value = "HV001"
$ws.Range("D2").Value = "HV002"
$ws.Range("D3").Value = "HV003"
$ws.Range("D4").Value = "HV004"
$ws.Range("D5").Value = "HV005"
$ws.Range("D6").Value = "HV001"
$ws.Range("D7").Value = "HV002"
$ws.Range("D8").Value = "HV003"
$ws.Range("D9").Value = "HV004"
$ws.Range("D10").Value = "HV005"

# Column E ("HP0X" -> "HP00X")
$ws.Range("E1").Value = "HP001"
$ws.Range("E2").Value = "HP002"
$ws.Range("E3").Value = "HP002"
$ws.Range("E4").Value = "HP005"
$ws.Range("E5").Value = "HP005"
$ws.Range("E6").Value = "HP002"
$ws.Range("E7").Value = "HP003"
$ws.Range("E8").Value = "HP004"
$ws.Range("E9").Value = "HP004"
$ws.Range("E10").Value = "HP004"

# Move the active selection from H6 to E10, as in the authored change.
$ws.Range("E10").Select() | Out-Null
